# Daily attendance processing - 2025-11-15 10:21:44
#
# For every row in column G ("Recorded By") whose value begins with the
# literal prefix "System, ", move that leading "System, " token to the
# end of the comma-separated list (as ", System") instead.
#
# Example:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"
#
# Cells whose value does not start with "System, " (e.g. plain
# "dnasr281@gmail.com", or the single word "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "System, "
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $newVal = $rest + ", System"
        $cell.Value2 = $newVal
        $changed = $changed + 1
    }
}

Write-Host "Updated Recorded By (column G) cells: $changed"
